# Slide 3 ("Built-in datasets") has an R-code content placeholder that
# demonstrates installing/loading the ggplot2 package. The deck was
# updated to install/load "ggplot2movies" instead (the sample uses the
# `movies` dataset that actually lives in that package).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# Find the content placeholder that holds the R code snippets by name
# (more robust than a hard-coded shape index).
$sh = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.Name -eq "Content Placeholder 2") {
        $sh = $cand
    }
}

$tr = $sh.TextFrame.TextRange
$paraCount = $tr.Paragraphs().Count

# Locate the two paragraphs we need to touch by their current text,
# instead of relying on fixed paragraph indices.
$installIdx = -1
$libraryIdx = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $t = $tr.Paragraphs($i, 1).Text
    if ($t -like '*install.packages("ggplot2")*') { $installIdx = $i }
    if ($t -like '*library(ggplot2)*') { $libraryIdx = $i }
}

# 1) > install.packages("ggplot2")  ->  > install.packages("ggplot2movies")
if ($installIdx -gt 0) {
    $installPara = $tr.Paragraphs($installIdx, 1)
    $installRange = $tr.Characters($installPara.Start, $installPara.Length)
    $installRange.Text = '> install.packages("ggplot2movies")'
}

# 2) > library(ggplot2)  ->  "> " (unchanged run) + "library(ggplot2movies)" (new run)
if ($libraryIdx -gt 0) {
    $libraryPara = $tr.Paragraphs($libraryIdx, 1)
    # Keep the leading "> " run untouched and retype only "library(ggplot2)"
    # so the call, retargeted at the new package, lands in its own run.
    $callRange = $tr.Characters($libraryPara.Start + 2, $libraryPara.Length - 2)
    $callRange.Text = "library(ggplot2movies)"
}
